# Automatische test-sync: 2025-08-03 23:36:50
# Adds a new log row (#54) to the "Logs" sheet for the new test mail about
# a refund ("Retour / Terugbetaling"), and updates the "Dashboard" sheet's
# category/count table to reflect the new totals (Retour / Terugbetaling
# goes from 2 -> 3 and moves up in the sorted-by-count list).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 54
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(54, 1).Value = "Ik heb nog geen geld terug."
$logs.Cells.Item(54, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(54, 3).Value = "Testmail #5: Ik heb nog geen geld terug."
$logs.Cells.Item(54, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(54, 5).Value = "Beste klant,`nHartelijk dank voor uw bericht. Kunt u mij het volgende laten weten: om welke aankoop gaat het precies en wanneer heeft u deze gedaan? Op basis van die informatie kunnen we de status van uw terugbetaling controleren en u verder helpen.`nMet vriendelijke groet,  `n[Naam]  `nKlantenservice Team"
$logs.Cells.Item(54, 6).Value = "2025-08-03 23:35:55"
$logs.Cells.Item(54, 7).Value = "Ja"
$logs.Cells.Item(54, 8).Value = "Nee"
$logs.Cells.Item(54, 9).Value = "Ja"
$logs.Cells.Item(54, 10).Value = "Nee"

# The multi-line text in column E would otherwise leave the row with an
# explicit (autofit) custom height; restore the default row height so the
# new row matches the rest of the sheet (no custom row height).
$logs.Rows.Item(54).EntireRow.AutoFit()

# Extend the conditional-formatting ranges (D, G, H, I, J) so they cover
# the newly added row 54 as well.
$oldRanges = @("D2:D53", "G2:G53", "H2:H53", "I2:I53", "J2:J53")
$newRanges = @("D2:D54", "G2:G54", "H2:H54", "I2:I54", "J2:J54")

for ($i = 0; $i -lt $oldRanges.Count; $i++) {
    $oldRange = $logs.Range($oldRanges[$i])
    $newRange = $logs.Range($newRanges[$i])
    $fc = $oldRange.FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($newRange)
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: update the Categorie / Aantal breakdown table
#    "Retour / Terugbetaling" count increases from 2 to 3 which moves it
#    up above the other categories that still have a count of 2.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(6, 2).Value = 3

$dash.Cells.Item(7, 1).Value = "Documentatie / Datasheets"
$dash.Cells.Item(7, 2).Value = 2

$dash.Cells.Item(8, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(8, 2).Value = 2

$dash.Cells.Item(9, 1).Value = "Klantenservice / Contact"
$dash.Cells.Item(9, 2).Value = 2
